$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commodities")

$ws.Range("B4").Value = 'MAT'
$ws.Range("C4").Value = 'iip_steel_oxygen'

$ws.Range("B5").Value = 'MAT'
$ws.Range("C5").Value = 'iip_steel_sinter'

$ws.Range("B6").Value = 'MAT'
$ws.Range("C6").Value = 'iip_steel_scrap'

$ws.Range("C7").Value = 'pri_uran'

$ws.Range("B8").Value = 'ENV'
$ws.Range("C8").Value = 'emi_CO2_f_x2x_neg_reusable'

$ws.Range("C9").Value = 'sec_heat_high'

$ws.Range("B10").Value = 'NRG'
$ws.Range("C10").Value = 'sec_heavy_fuel_oil'

$ws.Range("B11").Value = 'ENV'
$ws.Range("C11").Value = 'emi_CH4_f_ind'

$ws.Range("B12").Value = 'NRG'
$ws.Range("C12").Value = 'CO2_f_pow'

$ws.Range("C13").Value = 'pri_biomass'

$ws.Range("B14").Value = 'ENV'
$ws.Range("C14").Value = '[emi_CO2_f_x2x_neg_reusable'

$ws.Range("B15").Value = 'ENV'
$ws.Range("C15").Value = 'emi_CO2_f_ind'

$ws.Range("B16").Value = 'MAT'
$ws.Range("C16").Value = 'iip_steel_sponge_iron'

$ws.Range("B17").Value = 'NRG'
$ws.Range("C17").Value = 'iip_coke'

$ws.Range("B18").Value = 'MAT'
$ws.Range("C18").Value = 'iip_steel_crudesteel'

$ws.Range("C19").Value = 'sec_heat_low'

$ws.Range("C20").Value = 'pri_waste'

$ws.Range("B22").Value = 'NRG'
$ws.Range("C22").Value = 'pri_crude_oil'

$ws.Range("B23").Value = 'MAT'
$ws.Range("C23").Value = 'iip_steel_raw_iron'

$ws.Range("B24").Value = 'NRG'
$ws.Range("C24").Value = 'sec_biogas'

$ws.Range("B25").Value = 'NRG'
$ws.Range("C25").Value = 'pri_hydro_energy'

$ws.Range("C26").Value = 'sec_elec_ind'

$ws.Range("C27").Value = 'sec_natural_gas_syn'

$ws.Range("C28").Value = 'pri_geoth_heat'

$ws.Range("C29").Value = 'sec_H2'

$ws.Range("C30").Value = 'iip_heat_proc'

$ws.Range("C31").Value = 'iip_steel_iron_pellets'

$ws.Range("C32").Value = 'sec_elec'

$ws.Range("B33").Value = 'NRG'
$ws.Range("C33").Value = 'pri_natural_gas'

$ws.Range("B34").Value = 'ENV'
$ws.Range("C34").Value = 'emi_CO2_f_x2x_neg_stored]'

$ws.Range("B35").Value = 'DEM'
$ws.Range("C35").Value = 'exo_steel'

$ws.Range("B36").Value = 'ENV'
$ws.Range("C36").Value = 'emi_N2O_f_ind'

$ws.Range("B37").Value = 'NRG'
$ws.Range("C37").Value = 'pri_coal'

$ws.Range("B38").Value = 'NRG'
$ws.Range("C38").Value = 'sec_heating_oil'

$ws.Range("C39").Value = 'iip_steel_blafu_slag'

$ws.Range("C40").Value = 'pri_solar_radiation'

$ws.Range("B41").Value = 'NRG'
$ws.Range("C41").Value = 'sec_hydrogen'

$ws.Range("B42").Value = 'NRG'
$ws.Range("C42").Value = 'pri_wind_energy_on'

$ws.Range("B43").Value = 'MAT'
$ws.Range("C43").Value = 'iip_steel_iron_ore'

$ws.Range("B44").Value = 'NRG'
$ws.Range("C44").Value = 'sec_methane'
